# Auto-generated edit script: applies numeric cell updates across
# sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1636.5454
$ws.Range("I100").Value = 1500.2
$ws.Range("K100").Value = 1500.2
$ws.Range("M100").Value = -959.2
$ws.Range("H132").Value = 1466
$ws.Range("I132").Value = 1492.4348
$ws.Range("J132").Value = 1364.6666
$ws.Range("K132").Value = 4477.3044
$ws.Range("L132").Value = 4093.9998
$ws.Range("M132").Value = -1947.3044
$ws.Range("N132").Value = -9153.9998
$ws.Range("H137").Value = 4666.6665
$ws.Range("I137").Value = 3000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -6450
$ws.Range("H138").Value = 4188.6
$ws.Range("I138").Value = 2367.2
$ws.Range("K138").Value = 7101.599999999999
$ws.Range("M138").Value = -1961.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3164.12
$ws.Range("I2").Value = 2304.05
$ws.Range("K2").Value = 2304.05
$ws.Range("M2").Value = -2191.05
$ws.Range("H88").Value = 1441.5217
$ws.Range("I88").Value = 1455.1818
$ws.Range("J88").Value = 1429
$ws.Range("K88").Value = 1455.1818
$ws.Range("L88").Value = 1429
$ws.Range("M88").Value = -1049.1818
$ws.Range("N88").Value = -2241
$ws.Range("H91").Value = 1441.5217
$ws.Range("I91").Value = 1455.1818
$ws.Range("J91").Value = 1429
$ws.Range("K91").Value = 1455.1818
$ws.Range("L91").Value = 1429
$ws.Range("M91").Value = -51.18180000000007
$ws.Range("N91").Value = -4237
$ws.Range("H97").Value = 1385.25
$ws.Range("I97").Value = 1175.2667
$ws.Range("K97").Value = 1175.2667
$ws.Range("M97").Value = -679.2666999999999
$ws.Range("H102").Value = 2978.125
$ws.Range("I102").Value = 2025
$ws.Range("K102").Value = 2025
$ws.Range("M102").Value = -403
$ws.Range("H116").Value = 3164.12
$ws.Range("I116").Value = 2304.05
$ws.Range("K116").Value = 2304.05
$ws.Range("M116").Value = -10.05000000000018
$ws.Range("H122").Value = 2779.6155
$ws.Range("I122").Value = 2375.5454
$ws.Range("K122").Value = 7126.6362
$ws.Range("M122").Value = -4676.6362
$ws.Range("H132").Value = 1400.3414
$ws.Range("I132").Value = 1235.35
$ws.Range("K132").Value = 3706.05
$ws.Range("M132").Value = -1176.05
$ws.Range("H139").Value = 99749.5
$ws.Range("J139").Value = 99749.5
$ws.Range("L139").Value = 99749.5
$ws.Range("N139").Value = -110029.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3164.12
$ws.Range("I3").Value = 2304.05
$ws.Range("K3").Value = 2304.05
$ws.Range("M3").Value = -2190.05
$ws.Range("H94").Value = 1299.7059
$ws.Range("I94").Value = 1186.1538
$ws.Range("K94").Value = 1186.1538
$ws.Range("M94").Value = -735.1538
$ws.Range("H99").Value = 3928.4062
$ws.Range("I99").Value = 3032.5
$ws.Range("K99").Value = 3032.5
$ws.Range("M99").Value = -1534.5
$ws.Range("H107").Value = 1635.0769
$ws.Range("I107").Value = 1652.1111
$ws.Range("K107").Value = 1652.1111
$ws.Range("M107").Value = 267.8888999999999
$ws.Range("H134").Value = 3360.449
$ws.Range("I134").Value = 3170.0625
$ws.Range("K134").Value = 9510.1875
$ws.Range("M134").Value = -6975.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2986.842
$ws.Range("I16").Value = 2397.25
$ws.Range("J16").Value = 3997.5715
$ws.Range("K16").Value = 2397.25
$ws.Range("L16").Value = 3997.5715
$ws.Range("M16").Value = -2110.25
$ws.Range("N16").Value = -4571.5715
$ws.Range("H31").Value = 10627.692
$ws.Range("I31").Value = 18288.111
$ws.Range("J31").Value = 6572.1763
$ws.Range("K31").Value = 18288.111
$ws.Range("L31").Value = 6572.1763
$ws.Range("M31").Value = -17993.111
$ws.Range("N31").Value = -7162.1763
$ws.Range("H34").Value = 10627.692
$ws.Range("I34").Value = 18288.111
$ws.Range("J34").Value = 6572.1763
$ws.Range("K34").Value = 18288.111
$ws.Range("L34").Value = 6572.1763
$ws.Range("M34").Value = -18086.111
$ws.Range("N34").Value = -6976.1763
$ws.Range("H41").Value = 15479.154
$ws.Range("J41").Value = 15479.154
$ws.Range("L41").Value = 15479.154
$ws.Range("N41").Value = -16335.154
$ws.Range("H99").Value = 8363.799999999999
$ws.Range("I99").Value = 8447.5
$ws.Range("K99").Value = 8447.5
$ws.Range("M99").Value = -6949.5
$ws.Range("H102").Value = 360000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H103").Value = 4008
$ws.Range("I103").Value = 4008
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 4008
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -2836
$ws.Range("N103").ClearContents()
$ws.Range("H109").Value = 30142.5
$ws.Range("J109").Value = 30142.5
$ws.Range("L109").Value = 30142.5
$ws.Range("N109").Value = -32222.5
$ws.Range("H113").Value = 2986.842
$ws.Range("I113").Value = 2397.25
$ws.Range("J113").Value = 3997.5715
$ws.Range("K113").Value = 2397.25
$ws.Range("L113").Value = 3997.5715
$ws.Range("M113").Value = -227.25
$ws.Range("N113").Value = -8337.5715
$ws.Range("H126").Value = 8363.799999999999
$ws.Range("I126").Value = 8447.5
$ws.Range("K126").Value = 25342.5
$ws.Range("M126").Value = -22872.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 506.06668
$ws.Range("J98").Value = 283.8
$ws.Range("L98").Value = 851.4000000000001
$ws.Range("N98").Value = -3847.4
$ws.Range("H112").Value = 6306.5
$ws.Range("I112").Value = 6846.3335
$ws.Range("K112").Value = 20539.0005
$ws.Range("M112").Value = -19431.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 40200
$ws.Range("J93").Value = 40200
$ws.Range("L93").Value = 40200
$ws.Range("N93").Value = -43944
$ws.Range("H126").Value = 2946.7083
$ws.Range("I126").Value = 2157.75
$ws.Range("J126").Value = 4524.625
$ws.Range("K126").Value = 6473.25
$ws.Range("L126").Value = 13573.875
$ws.Range("M126").Value = -4003.25
$ws.Range("N126").Value = -18513.875
$ws.Range("H135").Value = 52036.668
$ws.Range("J135").Value = 52036.668
$ws.Range("L135").Value = 52036.668
$ws.Range("N135").Value = -62176.668
$ws.Range("H137").Value = 61750
$ws.Range("J137").Value = 61700
$ws.Range("L137").Value = 61700
$ws.Range("N137").Value = -71900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5094.5
$ws.Range("I7").Value = 5520
$ws.Range("K7").Value = 5520
$ws.Range("M7").Value = -5408
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1093
$ws.Range("H61").Value = 94130.73
$ws.Range("I61").Value = 114660.336
$ws.Range("J61").Value = 1747.5
$ws.Range("K61").Value = 114660.336
$ws.Range("L61").Value = 1747.5
$ws.Range("M61").Value = -114458.336
$ws.Range("N61").Value = -2151.5
$ws.Range("H113").Value = 94130.73
$ws.Range("I113").Value = 114660.336
$ws.Range("J113").Value = 1747.5
$ws.Range("K113").Value = 114660.336
$ws.Range("L113").Value = 1747.5
$ws.Range("M113").Value = -112490.336
$ws.Range("N113").Value = -6087.5
$ws.Range("H126").Value = 5094.5
$ws.Range("I126").Value = 5520
$ws.Range("K126").Value = 16560
$ws.Range("M126").Value = -14090
$ws.Range("H136").Value = 5854.6294
$ws.Range("I136").Value = 5772.115
$ws.Range("K136").Value = 17316.345
$ws.Range("M136").Value = -14766.345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10366.143
$ws.Range("I45").Value = 8855
$ws.Range("J45").Value = 11499.5
$ws.Range("K45").Value = 8855
$ws.Range("L45").Value = 11499.5
$ws.Range("M45").Value = -8364
$ws.Range("N45").Value = -12481.5
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77774
$ws.Range("H132").Value = 4168.42
$ws.Range("I132").Value = 3513.244
$ws.Range("J132").Value = 7153.1113
$ws.Range("K132").Value = 10539.732
$ws.Range("L132").Value = 21459.3339
$ws.Range("M132").Value = -8009.732
$ws.Range("N132").Value = -26519.3339
